$wb = $excel.ActiveWorkbook

# --- Active sheet / selection changes -------------------------------------------------
# Before: "entities" (index 1) was the active/selected tab with selection D15.
# After : "attributes" (index 2) becomes the active/selected tab with selection E17;
#         "entities" keeps its own view but is no longer the selected tab.
$wsEntities = $wb.Worksheets.Item("entities")
$wsEntities.Activate()
$wsEntities.Range("D15").Select() | Out-Null

$wsAttributes = $wb.Worksheets.Item("attributes")
$wsAttributes.Activate()
$wsAttributes.Range("E17").Select() | Out-Null

# --- Column width tweaks on the "attributes" sheet -------------------------------------
# (Excel's ColumnWidth is in "characters"; the stored OOXML width is ColumnWidth + 5/6.)
$wsAttributes.Columns.Item(2).ColumnWidth = 15.96 - 0.8333333333333334
$wsAttributes.Columns.Item(3).ColumnWidth = 23.07 - 0.8333333333333334
$wsAttributes.Columns.Item(6).ColumnWidth = 13.35 - 0.8333333333333334
$wsAttributes.Columns.Item(7).ColumnWidth = 12.9 - 0.8333333333333334

# --- Rebuild the "unique" column (H) and turn "nillable" (G) into TRUE()/FALSE() formulas ---
$rows = 2,3,4,5,6,7,8,9,10,11,12,13
foreach ($r in $rows) {
    $gCell = $wsAttributes.Range("G" + $r)
    if ($r -eq 2) {
        $gCell.Formula = "=TRUE()"
    } else {
        $gCell.Formula = "=FALSE()"
    }

    $hCell = $wsAttributes.Range("H" + $r)
    $hCell.Value = $false
    $hCell.NumberFormat = $gCell.NumberFormat
}
